$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Стол N" shared-string labels in column A with plain numbers.
# New mapping per row (A2:A24):
$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 2
    8  = 2
    9  = 2
    10 = 2
    11 = 2
    12 = 3
    13 = 4
    14 = 4
    15 = 4
    16 = 4
    17 = 3
    18 = 5
    19 = 5
    20 = 5
    21 = 5
    22 = 4
    23 = 4
    24 = 4
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}

# Update the active selection to match the saved view state.
$ws.Range("B12").Select()
